$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.714.30"
$ws.Range("E2").Value = "  -4.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.99"
$ws.Range("E3").Value = "  -3.09%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "277.20"
$ws.Range("E5").Value = "  -8.04%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5083"
$ws.Range("E7").Value = "  -5.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3529"
$ws.Range("E8").Value = "  -5.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.46"
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06654"
$ws.Range("E10").Value = "  -7.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.06"
$ws.Range("E11").Value = "  -7.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8247"
$ws.Range("E12").Value = "  -7.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07873"
$ws.Range("E13").Value = "  -3.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.829.04"
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.069"
$ws.Range("E15").Value = "  -4.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.50"
$ws.Range("E16").Value = "  -6.51%  "
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.10"
$ws.Range("E18").Value = "  -5.04%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000008015"
$ws.Range("E20").Value = "  -6.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "25.757.48"
$ws.Range("E21").Value = "  -4.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.738"
$ws.Range("E22").Value = "  -5.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.995"
$ws.Range("E23").Value = "  -6.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.093"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.14"
$ws.Range("E25").Value = "  -2.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.212"
$ws.Range("E26").Value = "  -4.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.675"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.09"
$ws.Range("E28").Value = "  -5.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "109.54"
$ws.Range("E29").Value = "  -3.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.332"
$ws.Range("E30").Value = "  -8.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.231"
$ws.Range("E31").Value = "  -8.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08791"
$ws.Range("E32").Value = "  -4.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04869"
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7255"
$ws.Range("E34").Value = "  -11.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("E36").Value = "  -2.57%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.126"
$ws.Range("E37").Value = "  -2.93%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.360"
$ws.Range("E38").Value = "  -9.90%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01850"
$ws.Range("E39").Value = "  -5.28%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5162"
$ws.Range("E40").Value = "  -14.87%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9633"
$ws.Range("E41").Value = "  -9.94%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.216"
$ws.Range("E42").Value = "  -6.39%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "110.27"
$ws.Range("E43").Value = "  -3.94%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.022"
$ws.Range("E44").Value = "  -10.07%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4550"
$ws.Range("E46").Value = "  -10.98%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1363"
$ws.Range("E47").Value = "  -8.79%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "36.56"
$ws.Range("E48").Value = "  -3.21%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.265"
$ws.Range("E49").Value = "  -6.88%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.496"
$ws.Range("E50").Value = "  -8.40%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05841"
$ws.Range("E51").Value = "  -4.01%  "
